$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(308051846, Eyal  Sofer: -10,-2)"
$ws.Range("B1").Value = "(312049950, Molham  Peretz: -2,8)"
$ws.Range("C1").Value = "(308073899, Anan  Kirshenbaum: 9,-4)"
$ws.Range("D1").Value = "(318869187, Soaad  Leibovich: -1,-8)"
$ws.Range("E1").Value = "(205898513, Asaf  Braymok: -5,-9)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: -2,4)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: 2,7)"

$ws.Range("A3").Value = "cost: 483.3880298169999"
$ws.Range("A4").Value = "time: 91.6776059634"
